$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/3/2025  Through  11/9/2025"

# --- Cells changing from text placeholder ("N/A"/"0") to a real number ---
# Restore numeric formatting before assigning, so the cell reuses the standard
# "#,##0" numeric style instead of staying a text-formatted cell.
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 1
$ws.Range("F33").NumberFormat = "#,##0"
$ws.Range("F33").Value = 2

# --- Cells changing from a number to the text placeholder ("0" / "***.*") ---
# Copy formatting+value from an existing placeholder cell of the same kind so the
# destination picks up the shared "N/A" style instead of becoming a plain string.
$ws.Range("C14").Copy($ws.Range("F22"))
$ws.Range("C14").Copy($ws.Range("G31"))
$ws.Range("E14").Copy($ws.Range("H31"))

# --- Plain numeric value updates ---
$ws.Range("M14").Value = 200
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -87.5
$ws.Range("I15").Value = 24
$ws.Range("J15").Value = 39
$ws.Range("K15").Value = -38.461538461538
$ws.Range("L15").Value = 9.090909090909
$ws.Range("M15").Value = -14.285714285714
$ws.Range("N15").Value = -17.241379310344
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -66.666666666666
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = -38.888888888888
$ws.Range("I16").Value = 273
$ws.Range("J16").Value = 376
$ws.Range("K16").Value = -27.393617021276
$ws.Range("L16").Value = -10.78431372549
$ws.Range("M16").Value = -9.602649006622
$ws.Range("N16").Value = -75.840707964601
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 54
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 17.391304347826
$ws.Range("I17").Value = 480
$ws.Range("J17").Value = 560
$ws.Range("K17").Value = -14.285714285714
$ws.Range("L17").Value = 15.662650602409
$ws.Range("M17").Value = 69.611307420494
$ws.Range("N17").Value = 28.686327077748
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 137
$ws.Range("J18").Value = 191
$ws.Range("K18").Value = -28.2722513089
$ws.Range("L18").Value = 10.483870967741
$ws.Range("M18").Value = -46.692607003891
$ws.Range("N18").Value = -91.945914168136
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 74
$ws.Range("H19").Value = -27.027027027027
$ws.Range("I19").Value = 629
$ws.Range("J19").Value = 836
$ws.Range("K19").Value = -24.760765550239
$ws.Range("L19").Value = -15.343203230148
$ws.Range("M19").Value = 46.279069767441
$ws.Range("N19").Value = -49.840510366826
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -40
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -34.482758620689
$ws.Range("I20").Value = 212
$ws.Range("J20").Value = 269
$ws.Range("K20").Value = -21.189591078066
$ws.Range("L20").Value = -28.135593220339
$ws.Range("M20").Value = 9.844559585492
$ws.Range("N20").Value = -89.117043121149
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -16.129032258064
$ws.Range("F21").Value = 159
$ws.Range("G21").Value = 202
$ws.Range("H21").Value = -21.287128712871
$ws.Range("I21").Value = 1761
$ws.Range("J21").Value = 2275
$ws.Range("K21").Value = -22.593406593406
$ws.Range("L21").Value = -7.65600419507
$ws.Range("M21").Value = 17.792642140468
$ws.Range("N21").Value = -72.697674418604
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -100
$ws.Range("J22").Value = 70
$ws.Range("K22").Value = -42.857142857142
$ws.Range("L22").Value = -48.051948051948
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = -4.807692307692
$ws.Range("I24").Value = 1225
$ws.Range("J24").Value = 1744
$ws.Range("K24").Value = -29.759174311926
$ws.Range("L24").Value = -28.27868852459
$ws.Range("M24").Value = 24.618514750763
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -62.5
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = -13.513513513513
$ws.Range("I25").Value = 405
$ws.Range("J25").Value = 896
$ws.Range("K25").Value = -54.799107142857
$ws.Range("L25").Value = -52.961672473867
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = -8.333333333333
$ws.Range("F26").Value = 76
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = -8.433734939759
$ws.Range("I26").Value = 907
$ws.Range("J26").Value = 1149
$ws.Range("K26").Value = -21.061792863359
$ws.Range("L26").Value = 8.233890214797
$ws.Range("M26").Value = 11.425061425061
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -88.888888888888
$ws.Range("I27").Value = 35
$ws.Range("J27").Value = 58
$ws.Range("K27").Value = -39.655172413793
$ws.Range("L27").Value = -14.634146341463
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -25
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 15
$ws.Range("H28").Value = -26.666666666666
$ws.Range("I28").Value = 103
$ws.Range("J28").Value = 130
$ws.Range("K28").Value = -20.76923076923
$ws.Range("L28").Value = -20.155038759689
$ws.Range("M29").Value = -87.5
$ws.Range("N29").Value = -98.148148148148
$ws.Range("M30").Value = -83.333333333333
$ws.Range("N30").Value = -98
$ws.Range("I33").Value = 3
$ws.Range("K33").Value = -40
$ws.Range("L33").Value = -57.142857142857
